$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory text for several rows
$ws.Range("H4").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H7").Value = "line graph(s)"
$ws.Range("H8").Value = "line graph(s)"
$ws.Range("H11").Value = "bar chart(s)"
$ws.Range("H14").Value = "line graph(s)"

# Remove the entire "is_viewed" column (column I)
$ws.Range("I1:I14").Delete()
